# Apply Week 16 game log + season totals update (Giants Team Data.xlsx)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "YDS": append this week's individual play-by-play yardage logs
# to the existing space-separated number strings.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 1 5 5 3 4 2 1 6 6 -2 -1 -4 10 3 7 1 6 5 0 5 0 2 2 4 4 2 7"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " -4 8 8 8 -2 7 6 8 2 9 11 5 2 0 4 5 9 6 4 8 2"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 11 0 3 15 5 0 13 4 1 1 2 5 5 2 3 9 1 4 4 6 6 6 4 5 3 3 3 6 0"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 6 2 3 23 13 46 1 9 1 6 39 4 3 16 10 5"

# ---------------------------------------------------------------------
# Sheet "OFF": season offensive totals, Home (row2) and Road (row3)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value2 = 350
$offWs.Range("E2").Value2 = 17
$offWs.Range("F2").Value2 = 128
$offWs.Range("G2").Value2 = 86
$offWs.Range("H2").Value2 = 15
$offWs.Range("I2").Value2 = 14
$offWs.Range("J2").Value2 = 41
$offWs.Range("N2").Value2 = 43
$offWs.Range("O2").Value2 = 47
$offWs.Range("P2").Value2 = 25

$offWs.Range("B3").Value2 = 18
$offWs.Range("C3").Value2 = 352
$offWs.Range("E3").Value2 = 64
$offWs.Range("F3").Value2 = 208
$offWs.Range("G3").Value2 = 55
$offWs.Range("H3").Value2 = 64
$offWs.Range("I3").Value2 = 128
$offWs.Range("J3").Value2 = 106
$offWs.Range("L3").Value2 = 577
$offWs.Range("M3").Value2 = 349
$offWs.Range("Q3").Value2 = 1005

# ---------------------------------------------------------------------
# Sheet "DEF": season defensive totals, Home (row2) and Road (row3)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value2 = 399
$defWs.Range("F2").Value2 = 132
$defWs.Range("G2").Value2 = 118
$defWs.Range("H2").Value2 = 13
$defWs.Range("J2").Value2 = 62
$defWs.Range("N2").Value2 = 31
$defWs.Range("O2").Value2 = 38

$defWs.Range("B3").Value2 = 23
$defWs.Range("C3").Value2 = 355
$defWs.Range("E3").Value2 = 60
$defWs.Range("F3").Value2 = 225
$defWs.Range("G3").Value2 = 64
$defWs.Range("H3").Value2 = 36
$defWs.Range("I3").Value2 = 118
$defWs.Range("J3").Value2 = 116
$defWs.Range("L3").Value2 = 572
$defWs.Range("M3").Value2 = 388
$defWs.Range("Q3").Value2 = 1041

# ---------------------------------------------------------------------
# Sheet "ST": special teams. Row2 season totals plus distance logs
# (strings in B4:B6 / D3:D5).
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value2 = 132
$stWs.Range("D2").Value2 = 128
$stWs.Range("F2").Value2 = 360
$stWs.Range("G2").Value2 = 343
$stWs.Range("N2").Value2 = 51
$stWs.Range("O2").Value2 = 32

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 65 51"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 20 10"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 17 21"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 47 48 28 38 33 33 31 36"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " -4 0 0 14 0 39 0 1"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 3 13 0"

# ---------------------------------------------------------------------
# Sheet "TURNS": turnovers, Road row (row3)
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value2 = 15
$turnsWs.Range("D3").Value2 = 16
$turnsWs.Range("E3").Value2 = 17

# ---------------------------------------------------------------------
# Sheet "PEN": penalties, row2
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value2 = 35
$penWs.Range("D2").Value2 = 12
